$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 2
$ws.Range("O2").Value = 1.29
$ws.Range("P2").Value = 3.75
$ws.Range("Y2").Value = 11
$ws.Range("AE2").Value = 11
$ws.Range("AI2").Value = 7.5

# Row 3
$ws.Range("G3").Value = 5.75
$ws.Range("H3").Value = 3.9
$ws.Range("I3").Value = 1.57
$ws.Range("Q3").Value = 1.8
$ws.Range("R3").Value = 2
$ws.Range("S3").Value = 2.75
$ws.Range("T3").Value = 1.4
$ws.Range("Z3").Value = 29
$ws.Range("AB3").Value = 51

# Row 9
$ws.Range("G9").Value = 1.27
$ws.Range("H9").Value = 6.25
$ws.Range("I9").Value = 8.5
$ws.Range("J9").Value = 1.67
$ws.Range("L9").Value = 8
$ws.Range("W9").Value = 1.83
$ws.Range("X9").Value = 1.83
$ws.Range("Y9").Value = 9.5
$ws.Range("AB9").Value = 8.5
$ws.Range("AO9").Value = 251
